$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet and update title cell to reflect new "through" date
$ws.Name = "Through 2022-04-20"
$ws.Range("I1").Value = "2022 (through 04-20)"

# Update the new data points for 2022
$ws.Range("I2").Value = 161   # January
$ws.Range("I5").Value = 87    # April
$ws.Range("I14").Value = 523  # Total
